# Auto-generated edit script: refresh market-price-derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets,
# mirroring a scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
# Row 48
$ws.Range("H48").Value = 3748.8333
$ws.Range("I48").Value = 2999.3333
$ws.Range("J48").Value = 4498.3335
$ws.Range("K48").Value = 8997.999899999999
$ws.Range("L48").Value = 13495.0005
$ws.Range("M48").Value = -8705.999899999999
$ws.Range("N48").Value = -14079.0005
# Row 51
$ws.Range("H51").Value = 2696.35
$ws.Range("I51").Value = 2423.5715
$ws.Range("J51").Value = 3332.8333
$ws.Range("K51").Value = 2423.5715
$ws.Range("L51").Value = 3332.8333
$ws.Range("M51").Value = -1939.5715
$ws.Range("N51").Value = -4300.8333
# Row 56
$ws.Range("H56").Value = 3748.8333
$ws.Range("I56").Value = 2999.3333
$ws.Range("J56").Value = 4498.3335
$ws.Range("K56").Value = 8997.999899999999
$ws.Range("L56").Value = 13495.0005
$ws.Range("M56").Value = -8463.999899999999
$ws.Range("N56").Value = -14563.0005
# Row 74
$ws.Range("H74").Value = 88960.69500000001
$ws.Range("I74").Value = 147999.14
$ws.Range("J74").Value = 20082.5
$ws.Range("K74").Value = 147999.14
$ws.Range("L74").Value = 20082.5
$ws.Range("M74").Value = -147063.14
$ws.Range("N74").Value = -21954.5
# Row 77
$ws.Range("H77").Value = 88960.69500000001
$ws.Range("I77").Value = 147999.14
$ws.Range("J77").Value = 20082.5
$ws.Range("K77").Value = 739995.7000000001
$ws.Range("L77").Value = 100412.5
$ws.Range("M77").Value = -735315.7000000001
$ws.Range("N77").Value = -109772.5
# Row 103
$ws.Range("H103").Value = 374.83334
# Row 137
$ws.Range("H137").Value = 1452.7
$ws.Range("I137").Value = 1445.4231
$ws.Range("J137").Value = 1500
$ws.Range("K137").Value = 4336.2693
$ws.Range("L137").Value = 4500
$ws.Range("M137").Value = -1786.2693
$ws.Range("N137").Value = -9600

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6300.067
$ws.Range("I32").Value = 3145.5818
$ws.Range("J32").Value = 40999.4
$ws.Range("K32").Value = 3145.5818
$ws.Range("L32").Value = 40999.4
$ws.Range("M32").Value = -2858.5818
$ws.Range("N32").Value = -41573.4
# Row 102
$ws.Range("H102").Value = 1677.375
$ws.Range("I102").Value = 1677.375
$ws.Range("K102").Value = 1677.375
$ws.Range("M102").Value = -55.375
# Row 122
$ws.Range("H122").Value = 1579.025
$ws.Range("I122").Value = 1492.579
$ws.Range("K122").Value = 4477.737
$ws.Range("M122").Value = -2027.737

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 1176.0312
$ws.Range("I105").Value = 1101.7916
$ws.Range("J105").Value = 1398.75
$ws.Range("K105").Value = 1101.7916
$ws.Range("L105").Value = 1398.75
$ws.Range("M105").Value = 645.2084
$ws.Range("N105").Value = -4892.75
# Row 134
$ws.Range("H134").Value = 1496.1971
$ws.Range("I134").Value = 1430.0725
$ws.Range("J134").Value = 3777.5
$ws.Range("K134").Value = 4290.2175
$ws.Range("L134").Value = 11332.5
$ws.Range("M134").Value = -1755.2175
$ws.Range("N134").Value = -16402.5

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 3500
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2876
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 3500
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -14380
$ws.Range("N65").ClearContents()
# Row 68
$ws.Range("H68").Value = 84600
$ws.Range("I68").Value = 84600
$ws.Range("K68").Value = 84600
$ws.Range("M68").Value = -83851
# Row 71
$ws.Range("H71").Value = 84600
$ws.Range("I71").Value = 84600
$ws.Range("K71").Value = 253800
$ws.Range("M71").Value = -250056

$ws = $wb.Worksheets.Item("CUL")
# Row 49
$ws.Range("H49").Value = 1075
$ws.Range("I49").Value = 433.33334
$ws.Range("K49").Value = 1300.00002
$ws.Range("M49").Value = -1144.00002
# Row 129
$ws.Range("H129").Value = 1404.75
$ws.Range("I129").Value = 899.7143
$ws.Range("J129").Value = 2583.1667
$ws.Range("K129").Value = 2699.1429
$ws.Range("L129").Value = 7749.500100000001
$ws.Range("M129").Value = 2300.8571
$ws.Range("N129").Value = -17749.5001
# Row 131
$ws.Range("H131").Value = 37282.105
$ws.Range("I131").Value = 91951.63
$ws.Range("J131").Value = 1907.7059
$ws.Range("K131").Value = 275854.89
$ws.Range("L131").Value = 5723.1177
$ws.Range("M131").Value = -270814.89
$ws.Range("N131").Value = -15803.1177

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 218.0625
$ws.Range("I2").Value = 152.38461
$ws.Range("K2").Value = 152.38461
$ws.Range("M2").Value = -39.38461000000001
# Row 102
$ws.Range("H102").Value = 38467050
$ws.Range("I102").Value = 1802.3334
$ws.Range("K102").Value = 1802.3334
$ws.Range("M102").Value = -180.3334

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 3553.5
$ws.Range("I68").Value = 3664.4
$ws.Range("J68").Value = 2999
$ws.Range("K68").Value = 3664.4
$ws.Range("L68").Value = 2999
$ws.Range("M68").Value = -2915.4
$ws.Range("N68").Value = -4497
# Row 71
$ws.Range("H71").Value = 3553.5
$ws.Range("I71").Value = 3664.4
$ws.Range("J71").Value = 2999
$ws.Range("K71").Value = 18322
$ws.Range("L71").Value = 14995
$ws.Range("M71").Value = -14578
$ws.Range("N71").Value = -22483
# Row 82
$ws.Range("H82").Value = 2120.8667
$ws.Range("I82").Value = 2400.7
$ws.Range("J82").Value = 1561.2
$ws.Range("K82").Value = 2400.7
$ws.Range("L82").Value = 1561.2
$ws.Range("M82").Value = -2039.7
$ws.Range("N82").Value = -2283.2
# Row 85
$ws.Range("H85").Value = 2120.8667
$ws.Range("I85").Value = 2400.7
$ws.Range("J85").Value = 1561.2
$ws.Range("K85").Value = 2400.7
$ws.Range("L85").Value = 1561.2
$ws.Range("M85").Value = -1152.7
$ws.Range("N85").Value = -4057.2
# Row 100
$ws.Range("H100").Value = 7721
$ws.Range("I100").Value = 2415.8333
$ws.Range("J100").Value = 18331.334
$ws.Range("K100").Value = 2415.8333
$ws.Range("L100").Value = 18331.334
$ws.Range("M100").Value = -1874.8333
$ws.Range("N100").Value = -19413.334

$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 39000
$ws.Range("I70").Value = 39000
$ws.Range("K70").Value = 39000
$ws.Range("M70").Value = -38685
# Row 73
$ws.Range("H73").Value = 39000
$ws.Range("I73").Value = 39000
$ws.Range("K73").Value = 39000
$ws.Range("M73").Value = -37908
# Row 132
$ws.Range("H132").Value = 214289.6
$ws.Range("I132").Value = 1540.5476
$ws.Range("J132").Value = 2001381.6
$ws.Range("K132").Value = 4621.642800000001
$ws.Range("L132").Value = 6004144.800000001
$ws.Range("M132").Value = -2091.642800000001
$ws.Range("N132").Value = -6009204.800000001
